$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "[-, -, -, 'MEC-2B-Metalografia']"

$ws.Range("B3").Value = "-"
$ws.Range("C3").Value = "[-, -, -, 'MEC-2B-Metalografia']"

$ws.Range("B4").Value = "-"

$ws.Range("B6").Value = "-"
$ws.Range("F6").Value = "[-, -, -, 'MEC-2B-Metalografia']"

$ws.Range("B7").Value = "-"

$ws.Range("C8").Value = "[-, -, 'MEC-2B-Metalografia', -]"
